# Update the Tnfsf12-Tnfrsf12a NATMI LR-pair sheet with the new TPM-derived
# expression numbers for the "ECs" cluster (ligand Tnfsf12, receptor Tnfrsf12a),
# and the values that are re-derived from them (specificity scores and edge
# weights), matching the refreshed natmiOut_TPM pipeline output.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 5.115481333333334
$ws.Range("H2").Value = 15.346444
$ws.Range("I2").Value = 0.1917470154127355
$ws.Range("J2").Value = 0.1917470154127354
$ws.Range("M2").Value = 1.593166333333333
$ws.Range("N2").Value = 4.779498999999999
$ws.Range("O2").Value = 0.08683039593887641
$ws.Range("P2").Value = 0.0868303959388764
$ws.Range("Q2").Value = 8.149812639061778
$ws.Range("R2").Value = 73.34831375155601
$ws.Range("S2").Value = 0.01664946926838566
$ws.Range("T2").Value = 0.01664946926838565

# Row 3
$ws.Range("G3").Value = 5.115481333333334
$ws.Range("H3").Value = 15.346444
$ws.Range("I3").Value = 0.1917470154127355
$ws.Range("J3").Value = 0.1917470154127354
$ws.Range("O3").Value = 0.5056861612921587
$ws.Range("P3").Value = 0.5056861612921587
$ws.Range("Q3").Value = 47.46318871560356
$ws.Range("R3").Value = 427.1686984404321
$ws.Range("S3").Value = 0.0969638121632946
$ws.Range("T3").Value = 0.09696381216329458

# Row 4
$ws.Range("G4").Value = 5.115481333333334
$ws.Range("H4").Value = 15.346444
$ws.Range("I4").Value = 0.1917470154127355
$ws.Range("J4").Value = 0.1917470154127354
$ws.Range("M4").Value = 7.476516666666666
$ws.Range("O4").Value = 0.407483442768965
$ws.Range("P4").Value = 0.4074834427689649
$ws.Range("Q4").Value = 38.24598144668889
$ws.Range("R4").Value = 344.2138330202
$ws.Range("S4").Value = 0.07813373398105525
$ws.Range("T4").Value = 0.07813373398105522

# Row 5
$ws.Range("I5").Value = 0.4487184033275903
$ws.Range("J5").Value = 0.4487184033275903
$ws.Range("M5").Value = 1.593166333333333
$ws.Range("N5").Value = 4.779498999999999
$ws.Range("O5").Value = 0.08683039593887641
$ws.Range("P5").Value = 0.0868303959388764
$ws.Range("Q5").Value = 19.07185312348767
$ws.Range("R5").Value = 171.646678111389
$ws.Range("S5").Value = 0.03896239662599511
$ws.Range("T5").Value = 0.0389623966259951

# Row 6
$ws.Range("I6").Value = 0.4487184033275903
$ws.Range("J6").Value = 0.4487184033275903
$ws.Range("O6").Value = 0.5056861612921587
$ws.Range("P6").Value = 0.5056861612921587
$ws.Range("S6").Value = 0.2269106868798758
$ws.Range("T6").Value = 0.2269106868798758

# Row 7
$ws.Range("I7").Value = 0.4487184033275903
$ws.Range("J7").Value = 0.4487184033275903
$ws.Range("M7").Value = 7.476516666666666
$ws.Range("O7").Value = 0.407483442768965
$ws.Range("P7").Value = 0.4074834427689649
$ws.Range("Q7").Value = 89.50165764778333
$ws.Range("R7").Value = 805.51491883005
$ws.Range("S7").Value = 0.1828453198217195
$ws.Range("T7").Value = 0.1828453198217195

# Row 8
$ws.Range("I8").Value = 0.3595345812596742
$ws.Range("J8").Value = 0.3595345812596742
$ws.Range("M8").Value = 1.593166333333333
$ws.Range("N8").Value = 4.779498999999999
$ws.Range("O8").Value = 0.08683039593887641
$ws.Range("P8").Value = 0.0868303959388764
$ws.Range("Q8").Value = 15.28127813735589
$ws.Range("R8").Value = 137.531503236203
$ws.Range("S8").Value = 0.03121853004449565
$ws.Range("T8").Value = 0.03121853004449564

# Row 9
$ws.Range("I9").Value = 0.3595345812596742
$ws.Range("J9").Value = 0.3595345812596742
$ws.Range("O9").Value = 0.5056861612921587
$ws.Range("P9").Value = 0.5056861612921587
$ws.Range("S9").Value = 0.1818116622489884
$ws.Range("T9").Value = 0.1818116622489884

# Row 10
$ws.Range("I10").Value = 0.3595345812596742
$ws.Range("J10").Value = 0.3595345812596742
$ws.Range("M10").Value = 7.476516666666666
$ws.Range("O10").Value = 0.407483442768965
$ws.Range("P10").Value = 0.4074834427689649
$ws.Range("Q10").Value = 71.71299586959444
$ws.Range("R10").Value = 645.41696282635
$ws.Range("S10").Value = 0.1465043889661903
$ws.Range("T10").Value = 0.1465043889661903
